# Auto-generated Excel COM-interop script to apply market-data refresh diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 2206.25  # H62: 2397.1904 -> 2206.25
$ws.Cells.Item(62, 9).Value = 1924.8422  # I62: 2067.5334 -> 1924.8422
$ws.Cells.Item(62, 10).Value = 3275.6  # J62: 3221.3333 -> 3275.6
$ws.Cells.Item(62, 11).Value = 1924.8422  # K62: 2067.5334 -> 1924.8422
$ws.Cells.Item(62, 12).Value = 3275.6  # L62: 3221.3333 -> 3275.6
$ws.Cells.Item(62, 13).Value = -1300.8422  # M62: -1443.5334 -> -1300.8422
$ws.Cells.Item(62, 14).Value = -4523.6  # N62: -4469.3333 -> -4523.6

$ws.Cells.Item(65, 8).Value = 2206.25  # H65: 2397.1904 -> 2206.25
$ws.Cells.Item(65, 9).Value = 1924.8422  # I65: 2067.5334 -> 1924.8422
$ws.Cells.Item(65, 10).Value = 3275.6  # J65: 3221.3333 -> 3275.6
$ws.Cells.Item(65, 11).Value = 9624.210999999999  # K65: 10337.667 -> 9624.210999999999
$ws.Cells.Item(65, 12).Value = 16378  # L65: 16106.6665 -> 16378
$ws.Cells.Item(65, 13).Value = -6504.210999999999  # M65: -7217.666999999999 -> -6504.210999999999
$ws.Cells.Item(65, 14).Value = -22618  # N65: -22346.6665 -> -22618

$ws.Cells.Item(98, 8).Value = 1420.7587  # H98: 1511.1852 -> 1420.7587
$ws.Cells.Item(98, 10).Value = 1314  # J98: 1685.3334 -> 1314
$ws.Cells.Item(98, 12).Value = 1314  # L98: 1685.3334 -> 1314
$ws.Cells.Item(98, 14).Value = -4310  # N98: -4681.3334 -> -4310

$ws.Cells.Item(111, 8).Value = 3766.6667  # H111: 809.0909 -> 3766.6667
$ws.Cells.Item(111, 9).Value = 5300  # I111: 810 -> 5300
$ws.Cells.Item(111, 10).Value = 700  # J111: 800 -> 700
$ws.Cells.Item(111, 11).Value = 15900  # K111: 2430 -> 15900
$ws.Cells.Item(111, 12).Value = 2100  # L111: 2400 -> 2100
$ws.Cells.Item(111, 13).Value = -12833  # M111: 637 -> -12833
$ws.Cells.Item(111, 14).Value = -8234  # N111: -8534 -> -8234

$ws.Cells.Item(122, 8).Value = 1420.7587  # H122: 1511.1852 -> 1420.7587
$ws.Cells.Item(122, 10).Value = 1314  # J122: 1685.3334 -> 1314
$ws.Cells.Item(122, 12).Value = 3942  # L122: 5056.0002 -> 3942
$ws.Cells.Item(122, 14).Value = -8842  # N122: -9956.0002 -> -8842

$ws.Cells.Item(125, 8).Value = 443.55  # H125: 431.70834 -> 443.55
$ws.Cells.Item(125, 9).Value = 420.92856  # I125: 415.1875 -> 420.92856
$ws.Cells.Item(125, 10).Value = 496.33334  # J125: 464.75 -> 496.33334
$ws.Cells.Item(125, 11).Value = 3788.35704  # K125: 3736.6875 -> 3788.35704
$ws.Cells.Item(125, 12).Value = 4467.00006  # L125: 4182.75 -> 4467.00006
$ws.Cells.Item(125, 13).Value = -1328.35704  # M125: -1276.6875 -> -1328.35704
$ws.Cells.Item(125, 14).Value = -9387.00006  # N125: -9102.75 -> -9387.00006

$ws.Cells.Item(138, 8).Value = 3478835.8  # H138: 4174171.2 -> 3478835.8
$ws.Cells.Item(138, 9).Value = 6998.5  # I138: 5333.3335 -> 6998.5
$ws.Cells.Item(138, 10).Value = 3629785  # J138: 4512185 -> 3629785
$ws.Cells.Item(138, 11).Value = 20995.5  # K138: 16000.0005 -> 20995.5
$ws.Cells.Item(138, 12).Value = 10889355  # L138: 13536555 -> 10889355
$ws.Cells.Item(138, 13).Value = -15855.5  # M138: -10860.0005 -> -15855.5
$ws.Cells.Item(138, 14).Value = -10899635  # N138: -13546835 -> -10899635

$ws.Cells.Item(141, 8).Value = 2865  # H141: 3065 -> 2865
$ws.Cells.Item(141, 9).Value = 1797.5  # I141: 2095 -> 1797.5
$ws.Cells.Item(141, 10).Value = 5000  # J141: 3550 -> 5000
$ws.Cells.Item(141, 11).Value = 5392.5  # K141: 6285 -> 5392.5
$ws.Cells.Item(141, 12).Value = 15000  # L141: 10650 -> 15000
$ws.Cells.Item(141, 13).Value = -212.5  # M141: -1105 -> -212.5
$ws.Cells.Item(141, 14).Value = -25360  # N141: -21010 -> -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4463712.5  # H32: 5072322 -> 4463712.5
$ws.Cells.Item(32, 9).Value = 5309671.5  # I32: 6311383 -> 5309671.5
$ws.Cells.Item(32, 10).Value = 22429.584  # J32: 20765.77 -> 22429.584
$ws.Cells.Item(32, 11).Value = 5309671.5  # K32: 6311383 -> 5309671.5
$ws.Cells.Item(32, 12).Value = 22429.584  # L32: 20765.77 -> 22429.584
$ws.Cells.Item(32, 13).Value = -5309384.5  # M32: -6311096 -> -5309384.5
$ws.Cells.Item(32, 14).Value = -23003.584  # N32: -21339.77 -> -23003.584

$ws.Cells.Item(74, 8).Value = 7001065.5  # H74: 7412857.5 -> 7001065.5
$ws.Cells.Item(74, 9).Value = 11409836  # I74: 12550761 -> 11409836
$ws.Cells.Item(74, 11).Value = 11409836  # K74: 12550761 -> 11409836
$ws.Cells.Item(74, 13).Value = -11408962  # M74: -12549887 -> -11408962

$ws.Cells.Item(77, 8).Value = 7001065.5  # H77: 7412857.5 -> 7001065.5
$ws.Cells.Item(77, 9).Value = 11409836  # I77: 12550761 -> 11409836
$ws.Cells.Item(77, 11).Value = 57049180  # K77: 62753805 -> 57049180
$ws.Cells.Item(77, 13).Value = -57044812  # M77: -62749437 -> -57044812

$ws.Cells.Item(102, 8).Value = 35716160  # H102: 23811408 -> 35716160
$ws.Cells.Item(102, 9).Value = 35716160  # I102: 23811408 -> 35716160
$ws.Cells.Item(102, 11).Value = 35716160  # K102: 23811408 -> 35716160
$ws.Cells.Item(102, 13).Value = -35714538  # M102: -23809786 -> -35714538

$ws.Cells.Item(110, 8).Value = 1112845.6  # H110: 5005005.5 -> 1112845.6
$ws.Cells.Item(110, 9).Value = 1112845.6  # I110: 5005005.5 -> 1112845.6
$ws.Cells.Item(110, 11).Value = 1112845.6  # K110: 5005005.5 -> 1112845.6
$ws.Cells.Item(110, 13).Value = -1110800.6  # M110: -5002960.5 -> -1110800.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 414.94116  # H80: 590.0714 -> 414.94116
$ws.Cells.Item(80, 9).Value = 75.333336  # I80: 82 -> 75.333336
$ws.Cells.Item(80, 10).Value = 600.1818  # J80: 793.3 -> 600.1818
$ws.Cells.Item(80, 11).Value = 75.333336  # K80: 82 -> 75.333336
$ws.Cells.Item(80, 12).Value = 600.1818  # L80: 793.3 -> 600.1818
$ws.Cells.Item(80, 13).Value = 922.666664  # M80: 916 -> 922.666664
$ws.Cells.Item(80, 14).Value = -2596.1818  # N80: -2789.3 -> -2596.1818

$ws.Cells.Item(83, 8).Value = 414.94116  # H83: 590.0714 -> 414.94116
$ws.Cells.Item(83, 9).Value = 75.333336  # I83: 82 -> 75.333336
$ws.Cells.Item(83, 10).Value = 600.1818  # J83: 793.3 -> 600.1818
$ws.Cells.Item(83, 11).Value = 376.66668  # K83: 410 -> 376.66668
$ws.Cells.Item(83, 12).Value = 3000.909  # L83: 3966.5 -> 3000.909
$ws.Cells.Item(83, 13).Value = 4615.33332  # M83: 4582 -> 4615.33332
$ws.Cells.Item(83, 14).Value = -12984.909  # N83: -13950.5 -> -12984.909

$ws.Cells.Item(107, 8).Value = 1865.3182  # H107: 2519.8 -> 1865.3182
$ws.Cells.Item(107, 9).Value = 1604.1666  # I107: 1899.75 -> 1604.1666
$ws.Cells.Item(107, 10).Value = 2178.7  # J107: 2933.1667 -> 2178.7
$ws.Cells.Item(107, 11).Value = 1604.1666  # K107: 1899.75 -> 1604.1666
$ws.Cells.Item(107, 12).Value = 2178.7  # L107: 2933.1667 -> 2178.7
$ws.Cells.Item(107, 13).Value = 315.8334  # M107: 20.25 -> 315.8334
$ws.Cells.Item(107, 14).Value = -6018.7  # N107: -6773.1667 -> -6018.7

$ws.Cells.Item(134, 8).Value = 5202  # H134: 3990 -> 5202
$ws.Cells.Item(134, 9).Value = 5637.2  # I134: 3778.6 -> 5637.2
$ws.Cells.Item(134, 10).Value = 4476.6665  # J134: 4577.222 -> 4476.6665
$ws.Cells.Item(134, 11).Value = 16911.6  # K134: 11335.8 -> 16911.6
$ws.Cells.Item(134, 12).Value = 13429.9995  # L134: 13731.666 -> 13429.9995
$ws.Cells.Item(134, 13).Value = -14376.6  # M134: -8800.799999999999 -> -14376.6
$ws.Cells.Item(134, 14).Value = -18499.9995  # N134: -18801.666 -> -18499.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1475.1538  # H31: 1618.3857 -> 1475.1538
$ws.Cells.Item(31, 9).Value = 1035.5264  # I31: 1533.3334 -> 1035.5264
$ws.Cells.Item(31, 10).Value = 1616.7288  # J31: 1635.9828 -> 1616.7288
$ws.Cells.Item(31, 11).Value = 1035.5264  # K31: 1533.3334 -> 1035.5264
$ws.Cells.Item(31, 12).Value = 1616.7288  # L31: 1635.9828 -> 1616.7288
$ws.Cells.Item(31, 13).Value = -740.5264  # M31: -1238.3334 -> -740.5264
$ws.Cells.Item(31, 14).Value = -2206.7288  # N31: -2225.9828 -> -2206.7288

$ws.Cells.Item(34, 8).Value = 1475.1538  # H34: 1618.3857 -> 1475.1538
$ws.Cells.Item(34, 9).Value = 1035.5264  # I34: 1533.3334 -> 1035.5264
$ws.Cells.Item(34, 10).Value = 1616.7288  # J34: 1635.9828 -> 1616.7288
$ws.Cells.Item(34, 11).Value = 1035.5264  # K34: 1533.3334 -> 1035.5264
$ws.Cells.Item(34, 12).Value = 1616.7288  # L34: 1635.9828 -> 1616.7288
$ws.Cells.Item(34, 13).Value = -833.5264  # M34: -1331.3334 -> -833.5264
$ws.Cells.Item(34, 14).Value = -2020.7288  # N34: -2039.9828 -> -2020.7288

$ws.Cells.Item(58, 8).Value = 24288468  # H58: 24880850 -> 24288468
$ws.Cells.Item(58, 9).Value = 35175750  # I58: 32906408 -> 35175750
$ws.Cells.Item(58, 10).Value = 1445.9231  # J58: 1619.8 -> 1445.9231
$ws.Cells.Item(58, 11).Value = 35175750  # K58: 32906408 -> 35175750
$ws.Cells.Item(58, 12).Value = 1445.9231  # L58: 1619.8 -> 1445.9231
$ws.Cells.Item(58, 13).Value = -35175547  # M58: -32906205 -> -35175547
$ws.Cells.Item(58, 14).Value = -1851.9231  # N58: -2025.8 -> -1851.9231

$ws.Cells.Item(136, 8).Value = 24288468  # H136: 24880850 -> 24288468
$ws.Cells.Item(136, 9).Value = 35175750  # I136: 32906408 -> 35175750
$ws.Cells.Item(136, 10).Value = 1445.9231  # J136: 1619.8 -> 1445.9231
$ws.Cells.Item(136, 11).Value = 105527250  # K136: 98719224 -> 105527250
$ws.Cells.Item(136, 12).Value = 4337.7693  # L136: 4859.4 -> 4337.7693
$ws.Cells.Item(136, 13).Value = -105524700  # M136: -98716674 -> -105524700
$ws.Cells.Item(136, 14).Value = -9437.7693  # N136: -9959.4 -> -9437.7693

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(63, 8).Value = 3461.1  # H63: 4698.4 -> 3461.1
$ws.Cells.Item(63, 9).Value = 2342.2  # I63: 4592 -> 2342.2
$ws.Cells.Item(63, 10).Value = 4580  # J63: 4725 -> 4580
$ws.Cells.Item(63, 11).Value = 7026.599999999999  # K63: 13776 -> 7026.599999999999
$ws.Cells.Item(63, 12).Value = 13740  # L63: 14175 -> 13740
$ws.Cells.Item(63, 13).Value = -6277.599999999999  # M63: -13027 -> -6277.599999999999
$ws.Cells.Item(63, 14).Value = -15238  # N63: -15673 -> -15238

$ws.Cells.Item(66, 8).Value = 3461.1  # H66: 4698.4 -> 3461.1
$ws.Cells.Item(66, 9).Value = 2342.2  # I66: 4592 -> 2342.2
$ws.Cells.Item(66, 10).Value = 4580  # J66: 4725 -> 4580
$ws.Cells.Item(66, 11).Value = 21079.8  # K66: 41328 -> 21079.8
$ws.Cells.Item(66, 12).Value = 41220  # L66: 42525 -> 41220
$ws.Cells.Item(66, 13).Value = -17335.8  # M66: -37584 -> -17335.8
$ws.Cells.Item(66, 14).Value = -48708  # N66: -50013 -> -48708

$ws.Cells.Item(107, 8).Value = 1202.679  # H107: 1252.5476 -> 1202.679
$ws.Cells.Item(107, 9).Value = 512.4872  # I107: 529.8919 -> 512.4872
$ws.Cells.Item(107, 10).Value = 1843.5714  # J107: 1821.4468 -> 1843.5714
$ws.Cells.Item(107, 11).Value = 1537.4616  # K107: 1589.6757 -> 1537.4616
$ws.Cells.Item(107, 12).Value = 5530.7142  # L107: 5464.3404 -> 5530.7142
$ws.Cells.Item(107, 13).Value = 382.5383999999999  # M107: 330.3243000000002 -> 382.5383999999999
$ws.Cells.Item(107, 14).Value = -9370.7142  # N107: -9304.340400000001 -> -9370.7142

$ws.Cells.Item(131, 8).Value = 842.72095  # H131: 819.7879 -> 842.72095
$ws.Cells.Item(131, 10).Value = 933.9143  # J131: 940.12 -> 933.9143
$ws.Cells.Item(131, 12).Value = 2801.7429  # L131: 2820.36 -> 2801.7429
$ws.Cells.Item(131, 14).Value = -12881.7429  # N131: -12900.36 -> -12881.7429

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 47877.66  # H132: 54004.742 -> 47877.66
$ws.Cells.Item(132, 9).Value = 29188.945  # I132: 33608.72 -> 29188.945
$ws.Cells.Item(132, 10).Value = 146660.86  # J132: 147243.72 -> 146660.86
$ws.Cells.Item(132, 11).Value = 87566.83499999999  # K132: 100826.16 -> 87566.83499999999
$ws.Cells.Item(132, 12).Value = 439982.58  # L132: 441731.16 -> 439982.58
$ws.Cells.Item(132, 13).Value = -85036.83499999999  # M132: -98296.16 -> -85036.83499999999
$ws.Cells.Item(132, 14).Value = -445042.58  # N132: -446791.16 -> -445042.58

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1840.2354  # H7: 2128.4666 -> 1840.2354
$ws.Cells.Item(7, 9).Value = 1799  # I7: 2032.9231 -> 1799
$ws.Cells.Item(7, 10).Value = 2500  # J7: 2749.5 -> 2500
$ws.Cells.Item(7, 11).Value = 1799  # K7: 2032.9231 -> 1799
$ws.Cells.Item(7, 12).Value = 2500  # L7: 2749.5 -> 2500
$ws.Cells.Item(7, 13).Value = -1687  # M7: -1920.9231 -> -1687
$ws.Cells.Item(7, 14).Value = -2724  # N7: -2973.5 -> -2724

$ws.Cells.Item(40, 8).Value = 4754.364  # H40: 4827.091 -> 4754.364
$ws.Cells.Item(40, 9).Value = 4588.8887  # I40: 4677.778 -> 4588.8887
$ws.Cells.Item(40, 11).Value = 4588.8887  # K40: 4677.778 -> 4588.8887
$ws.Cells.Item(40, 13).Value = -4452.8887  # M40: -4541.778 -> -4452.8887

$ws.Cells.Item(122, 8).Value = 3990.6667  # H122: 3716.4119 -> 3990.6667
$ws.Cells.Item(122, 9).Value = 3320  # I122: 3323.3333 -> 3320
$ws.Cells.Item(122, 10).Value = 4469.7144  # J122: 4659.8 -> 4469.7144
$ws.Cells.Item(122, 11).Value = 9960  # K122: 9969.999899999999 -> 9960
$ws.Cells.Item(122, 12).Value = 13409.1432  # L122: 13979.4 -> 13409.1432
$ws.Cells.Item(122, 13).Value = -7510  # M122: -7519.999899999999 -> -7510
$ws.Cells.Item(122, 14).Value = -18309.1432  # N122: -18879.4 -> -18309.1432

$ws.Cells.Item(126, 8).Value = 1840.2354  # H126: 2128.4666 -> 1840.2354
$ws.Cells.Item(126, 9).Value = 1799  # I126: 2032.9231 -> 1799
$ws.Cells.Item(126, 10).Value = 2500  # J126: 2749.5 -> 2500
$ws.Cells.Item(126, 11).Value = 5397  # K126: 6098.7693 -> 5397
$ws.Cells.Item(126, 12).Value = 7500  # L126: 8248.5 -> 7500
$ws.Cells.Item(126, 13).Value = -2927  # M126: -3628.7693 -> -2927
$ws.Cells.Item(126, 14).Value = -12440  # N126: -13188.5 -> -12440

$ws.Cells.Item(136, 8).Value = 105881.96  # H136: 70795.336 -> 105881.96
$ws.Cells.Item(136, 9).Value = 52917.05  # I136: 34655.453 -> 52917.05
$ws.Cells.Item(136, 10).Value = 238294.25  # J136: 172644.1 -> 238294.25
$ws.Cells.Item(136, 11).Value = 158751.15  # K136: 103966.359 -> 158751.15
$ws.Cells.Item(136, 12).Value = 714882.75  # L136: 517932.3 -> 714882.75
$ws.Cells.Item(136, 13).Value = -156201.15  # M136: -101416.359 -> -156201.15
$ws.Cells.Item(136, 14).Value = -719982.75  # N136: -523032.3 -> -719982.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 8600  # H5: 0 -> 8600
$ws.Cells.Item(5, 10).Value = 8600  # J5: 0 -> 8600
$ws.Cells.Item(5, 12).Value = 8600  # L5: 0 -> 8600
$ws.Cells.Item(5, 14).Value = -8824  # N5: new cell

$ws.Cells.Item(113, 8).Value = 1874.0667  # H113: 2139.3076 -> 1874.0667
$ws.Cells.Item(113, 9).Value = 612.3  # I113: 727.875 -> 612.3
$ws.Cells.Item(113, 11).Value = 1836.9  # K113: 2183.625 -> 1836.9
$ws.Cells.Item(113, 13).Value = 333.1000000000001  # M113: -13.625 -> 333.1000000000001

$ws.Cells.Item(126, 8).Value = 2716.6667  # H126: 1900.0769 -> 2716.6667
$ws.Cells.Item(126, 9).Value = 1266.6666  # I126: 800.1111 -> 1266.6666
$ws.Cells.Item(126, 10).Value = 4166.6665  # J126: 4375 -> 4166.6665
$ws.Cells.Item(126, 11).Value = 3799.9998  # K126: 2400.3333 -> 3799.9998
$ws.Cells.Item(126, 12).Value = 12499.9995  # L126: 13125 -> 12499.9995
$ws.Cells.Item(126, 13).Value = -1329.9998  # M126: 69.66670000000022 -> -1329.9998
$ws.Cells.Item(126, 14).Value = -17439.9995  # N126: -18065 -> -17439.9995

$ws.Cells.Item(136, 8).Value = 160362.61  # H136: 208331.4 -> 160362.61
$ws.Cells.Item(136, 9).Value = 149600.58  # I136: 209200.8 -> 149600.58
$ws.Cells.Item(136, 10).Value = 172918.33  # J136: 207462 -> 172918.33
$ws.Cells.Item(136, 11).Value = 448801.74  # K136: 627602.3999999999 -> 448801.74
$ws.Cells.Item(136, 12).Value = 518754.99  # L136: 622386 -> 518754.99
$ws.Cells.Item(136, 13).Value = -446251.74  # M136: -625052.3999999999 -> -446251.74
$ws.Cells.Item(136, 14).Value = -523854.99  # N136: -627486 -> -523854.99
